# excel_import_path : finalise function + examples files update
#
# The example workbook shipped with the package is touched up by hand in
# Excel:
#   - the first sheet ("importr1") is renamed to "impexp1" to match the
#     renamed import/export helper, and
#   - that same (now first) sheet is made the active/selected tab when the
#     workbook is reopened, instead of "Autre onglet 1".

$wb = $excel.ActiveWorkbook

# Rename the first worksheet to match the new function name.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "impexp1"

# Make it the active tab (was "Autre onglet 1" before).
$ws1.Activate()
